# Insert 3 new data rows before the current row 589 (shifts 589-671 down to
# 592-674, which matches the target dataset exactly since the target is the
# same data shifted down by 3 with 3 brand-new rows at the top).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A589:A591").EntireRow.Insert()

# Shared/common column values for this data subset (identical across all
# rows for this sheet: Mercado/Region/Categoria/etc.).
$commonA = 10
$commonB = "Vega Modelo de Temuco"
$commonC = "La Araucanía"
$commonE = 9
$commonF = 100112008
$commonG = "Coliflor"
$commonH = "Sin especificar"
$commonI = "Primera"
$commonN = "`$/unidad"
$commonQ = 1
$commonR = "Hortaliza"

function Set-NewRow {
    param(
        [int]$Row,
        [double]$D,
        [double]$J,
        [double]$K,
        [double]$L,
        [double]$M,
        [string]$O,
        [double]$P
    )

    $ws.Cells.Item($Row, 1).Value = $commonA
    $ws.Cells.Item($Row, 2).Value = $commonB
    $ws.Cells.Item($Row, 3).Value = $commonC
    $ws.Cells.Item($Row, 4).Value = $D
    $ws.Cells.Item($Row, 5).Value = $commonE
    $ws.Cells.Item($Row, 6).Value = $commonF
    $ws.Cells.Item($Row, 7).Value = $commonG
    $ws.Cells.Item($Row, 8).Value = $commonH
    $ws.Cells.Item($Row, 9).Value = $commonI
    $ws.Cells.Item($Row, 10).Value = $J
    $ws.Cells.Item($Row, 11).Value = $K
    $ws.Cells.Item($Row, 12).Value = $L
    $ws.Cells.Item($Row, 13).Value = $M
    $ws.Cells.Item($Row, 14).Value = $commonN
    $ws.Cells.Item($Row, 15).Value = $O
    $ws.Cells.Item($Row, 16).Value = $P
    $ws.Cells.Item($Row, 17).Value = $commonQ
    $ws.Cells.Item($Row, 18).Value = $commonR
}

Set-NewRow 589 45127 3800 950  1000 974  "Región Metropolitana" 974
Set-NewRow 590 45127 1000 1200 1200 1200 "Región de O'Higgins"  1200
Set-NewRow 591 45127 3000 1000 1000 1000 "Región del Maule"     1000
